$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 31   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/15/2024  Through  4/21/2024"

# --- Weekly crime statistics table updates (rows 14-30) ---
# Row 14
$ws.Range("D14").NumberFormat = 'General'
$ws.Range("D14").Value = "'0"
$ws.Range("E14").NumberFormat = 'General'
$ws.Range("E14").Value = "'***.*"

# Row 15
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 7
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 250
$ws.Range("M15").Value = -30
$ws.Range("N15").Value = -68.181818181818

# Row 16
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 5
$ws.Range("F16").Value = 21
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 44
$ws.Range("J16").Value = 63
$ws.Range("K16").Value = -30.15873015873
$ws.Range("L16").Value = 57.142857142857
$ws.Range("M16").Value = -45
$ws.Range("N16").Value = -85.852090032154

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("E17").Value = -22.222222222222
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 28
$ws.Range("H17").Value = 10.714285714285
$ws.Range("I17").Value = 150
$ws.Range("J17").Value = 113
$ws.Range("K17").Value = 32.743362831858
$ws.Range("L17").Value = 48.514851485148
$ws.Range("M17").Value = 56.25
$ws.Range("N17").Value = -28.22966507177

# Row 18
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -28.571428571428
$ws.Range("I18").Value = 32
$ws.Range("J18").Value = 52
$ws.Range("K18").Value = -38.461538461538
$ws.Range("L18").Value = 28
$ws.Range("M18").Value = -54.285714285714
$ws.Range("N18").Value = -93.810444874274

# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -20
$ws.Range("G19").Value = 30
$ws.Range("H19").Value = -33.333333333333
$ws.Range("I19").Value = 107
$ws.Range("J19").Value = 129
$ws.Range("K19").Value = -17.054263565891
$ws.Range("L19").Value = -1.834862385321
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = -32.278481012658

# Row 20
$ws.Range("C20").NumberFormat = '#,##0'
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -12.5
$ws.Range("I20").Value = 31
$ws.Range("J20").Value = 35
$ws.Range("K20").Value = -11.428571428571
$ws.Range("L20").Value = 29.166666666666
$ws.Range("M20").Value = -39.215686274509
$ws.Range("N20").Value = -92.705882352941

# Row 21
$ws.Range("C21").Value = 27
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 91
$ws.Range("G21").Value = 99
$ws.Range("H21").Value = -8.080808080808
$ws.Range("I21").Value = 371
$ws.Range("J21").Value = 405
$ws.Range("K21").Value = -8.395061728395
$ws.Range("L21").Value = 27.931034482758
$ws.Range("M21").Value = -10.817307692307
$ws.Range("N21").Value = -77.487864077669

# Row 23
$ws.Range("C23").NumberFormat = '#,##0'
$ws.Range("C23").Value = 1
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 16.666666666666
$ws.Range("I23").Value = 25
$ws.Range("J23").Value = 30
$ws.Range("K23").Value = -16.666666666666
$ws.Range("L23").Value = 31.578947368421
$ws.Range("M23").Value = 108.333333333333

# Row 24
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = -27.272727272727
$ws.Range("F24").Value = 104
$ws.Range("G24").Value = 88
$ws.Range("H24").Value = 18.181818181818
$ws.Range("I24").Value = 422
$ws.Range("J24").Value = 361
$ws.Range("K24").Value = 16.897506925207
$ws.Range("L24").Value = 37.459283387622
$ws.Range("M24").Value = 9.610389610389

# Row 25
$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = -38.461538461538
$ws.Range("F25").Value = 47
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = 34.285714285714
$ws.Range("I25").Value = 201
$ws.Range("J25").Value = 132
$ws.Range("K25").Value = 52.272727272727
$ws.Range("L25").Value = 97.058823529411

# Row 26
$ws.Range("C26").Value = 21
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = 61.538461538461
$ws.Range("F26").Value = 64
$ws.Range("G26").Value = 60
$ws.Range("H26").Value = 6.666666666666
$ws.Range("I26").Value = 202
$ws.Range("J26").Value = 211
$ws.Range("K26").Value = -4.265402843601
$ws.Range("L26").Value = 4.663212435233
$ws.Range("M26").Value = -43.098591549295

# Row 27
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 14
$ws.Range("K27").Value = 55.555555555555
$ws.Range("L27").Value = 100

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("D28").Value = 3
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E28").Value = -66.666666666666
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = 20
$ws.Range("I28").Value = 28
$ws.Range("J28").Value = 20
$ws.Range("K28").Value = 40
$ws.Range("L28").Value = 55.555555555555

# Row 29
$ws.Range("D29").NumberFormat = 'General'
$ws.Range("D29").Value = "'0"
$ws.Range("E29").NumberFormat = 'General'
$ws.Range("E29").Value = "'***.*"

# Row 30
$ws.Range("D30").NumberFormat = 'General'
$ws.Range("D30").Value = "'0"
$ws.Range("E30").NumberFormat = 'General'
$ws.Range("E30").Value = "'***.*"

# --- Column E width (bestfit auto-narrow after data update) ---
$ws.Range("E1").ColumnWidth = $ws.Range("D1").ColumnWidth